$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 19: new time-log entry ---
$ws.Range("A19").Value = 44460
$ws.Range("B19").Value = 0.85763888888888884
$ws.Range("C19").Value = 0.9458333333333333

# Match the duration-format style used by the existing D column entries
# (D18 carries the custom [h]:mm:ss number format) before writing the
# shared "end - start" formula so D19 inherits the same look.
$ws.Range("D19").NumberFormat = $ws.Range("D18").NumberFormat
$ws.Range("D19").Formula = "=C19-B19"

$ws.Range("E19").Value = "Project Plan Document"
$ws.Range("F19").Value = "COCOMO Model for Project Plan."

# --- Move the active selection the way the author left it ---
$ws.Range("E20").Select()
